$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 651
$ws.Range("B651").Value = 5159852
$ws.Range("F651").Value = 'Vitesse'
$ws.Range("G651").Value = 'FC Twente'
$ws.Range("H651").Value = 2
$ws.Range("I651").Value = 2
$ws.Range("J651").Value = 'D'
$ws.Range("K651").Value = 3.3
$ws.Range("M651").Value = 2.15
$ws.Range("N651").Value = 4
$ws.Range("O651").Value = 3.4
$ws.Range("P651").Value = 2
$ws.Range("R651").Value = 1.825
$ws.Range("S651").Value = 2.025
$ws.Range("T651").Value = 2.25
$ws.Range("U651").Value = 1.825
$ws.Range("V651").Value = 2.025
$ws.Range("X651").Value = 2.4
$ws.Range("Y651").Value = -1
$ws.Range("Z651").Value = 0.825
$ws.Range("AA651").Value = -1
$ws.Range("AB651").Value = 0.825

# Row 652
$ws.Range("B652").Value = 5159851
$ws.Range("F652").Value = 'Go Ahead Eagles'
$ws.Range("G652").Value = 'AZ'
$ws.Range("H652").Value = 1
$ws.Range("I652").Value = 4
$ws.Range("J652").Value = 'A'
$ws.Range("K652").Value = 3.75
$ws.Range("M652").Value = 2
$ws.Range("N652").Value = 4.2
$ws.Range("O652").Value = 3.6
$ws.Range("P652").Value = 1.85
$ws.Range("R652").Value = 2.025
$ws.Range("S652").Value = 1.825
$ws.Range("T652").Value = 2.75
$ws.Range("U652").Value = 1.975
$ws.Range("V652").Value = 1.875
$ws.Range("X652").Value = -1
$ws.Range("Y652").Value = 0.8500000000000001
$ws.Range("Z652").Value = -1
$ws.Range("AA652").Value = 0.825
$ws.Range("AB652").Value = 0.9750000000000001

# Row 654
$ws.Range("B654").Value = 5159591
$ws.Range("F654").Value = 'FC Utrecht'
$ws.Range("G654").Value = 'Excelsior'
$ws.Range("H654").Value = 1
$ws.Range("K654").Value = 1.45
$ws.Range("L654").Value = 4.5
$ws.Range("M654").Value = 7
$ws.Range("N654").Value = 1.4
$ws.Range("O654").Value = 4.75
$ws.Range("P654").Value = 8
$ws.Range("Q654").Value = -1.25
$ws.Range("R654").Value = 1.86
$ws.Range("S654").Value = 2.04
$ws.Range("T654").Value = 3.25
$ws.Range("U654").Value = 2.05
$ws.Range("V654").Value = 1.8
$ws.Range("W654").Value = 0.3999999999999999
$ws.Range("Z654").Value = -0.5
$ws.Range("AA654").Value = 0.52
$ws.Range("AC654").Value = 0.8

# Row 655
$ws.Range("B655").Value = 5159854
$ws.Range("F655").Value = 'Feyenoord'
$ws.Range("G655").Value = 'NEC'
$ws.Range("H655").Value = 2
$ws.Range("K655").Value = 1.285
$ws.Range("L655").Value = 5.75
$ws.Range("M655").Value = 10
$ws.Range("N655").Value = 1.285
$ws.Range("O655").Value = 5.25
$ws.Range("P655").Value = 11
$ws.Range("Q655").Value = -1.5
$ws.Range("R655").Value = 1.875
$ws.Range("S655").Value = 1.975
$ws.Range("T655").Value = 3
$ws.Range("U655").Value = 1.95
$ws.Range("V655").Value = 1.9
$ws.Range("W655").Value = 0.2849999999999999
$ws.Range("Z655").Value = 0.875
$ws.Range("AA655").Value = -1
$ws.Range("AC655").Value = 0.8999999999999999

# Row 755
$ws.Range("B755").Value = 5159401
$ws.Range("F755").Value = 'Ajax'
$ws.Range("G755").Value = 'FC Emmen'
$ws.Range("H755").Value = 3
$ws.Range("I755").Value = 1
$ws.Range("J755").Value = 'H'
$ws.Range("K755").Value = 1.125
$ws.Range("L755").Value = 8
$ws.Range("M755").Value = 21
$ws.Range("N755").Value = 1.111
$ws.Range("O755").Value = 9
$ws.Range("P755").Value = 21
$ws.Range("Q755").Value = -2.5
$ws.Range("R755").Value = 1.925
$ws.Range("S755").Value = 1.925
$ws.Range("T755").Value = 3.5
$ws.Range("U755").Value = 1.85
$ws.Range("V755").Value = 2
$ws.Range("W755").Value = 0.111
$ws.Range("Y755").Value = -1
$ws.Range("AA755").Value = 0.925
$ws.Range("AB755").Value = 0.8500000000000001
$ws.Range("AC755").Value = -1

# Row 756
$ws.Range("B756").Value = 6511816
$ws.Range("F756").Value = 'Fortuna Sittard'
$ws.Range("G756").Value = 'AZ'
$ws.Range("H756").Value = 0
$ws.Range("I756").Value = 3
$ws.Range("J756").Value = 'A'
$ws.Range("K756").Value = 5
$ws.Range("L756").Value = 3.8
$ws.Range("M756").Value = 1.65
$ws.Range("N756").Value = 4.5
$ws.Range("O756").Value = 4
$ws.Range("P756").Value = 1.727
$ws.Range("Q756").Value = 0.75
$ws.Range("R756").Value = 1.95
$ws.Range("S756").Value = 1.9
$ws.Range("T756").Value = 2.75
$ws.Range("U756").Value = 2.025
$ws.Range("V756").Value = 1.825
$ws.Range("W756").Value = -1
$ws.Range("Y756").Value = 0.7270000000000001
$ws.Range("AA756").Value = 0.8999999999999999
$ws.Range("AB756").Value = 0.5125
$ws.Range("AC756").Value = -0.5

# Row 901
$ws.Range("B901").Value = 6838426
$ws.Range("F901").Value = 'Sparta Rotterdam'
$ws.Range("G901").Value = 'Almere City FC'
$ws.Range("H901").Value = 1
$ws.Range("I901").Value = 2
$ws.Range("J901").Value = 'A'
$ws.Range("K901").Value = 2.15
$ws.Range("L901").Value = 3.4
$ws.Range("M901").Value = 3
$ws.Range("N901").Value = 1.7
$ws.Range("O901").Value = 4
$ws.Range("P901").Value = 4.75
$ws.Range("Q901").Value = -0.75
$ws.Range("R901").Value = 1.9
$ws.Range("S901").Value = 1.95
$ws.Range("T901").Value = 2.75
$ws.Range("U901").Value = 2.025
$ws.Range("V901").Value = 1.825
$ws.Range("W901").Value = -1
$ws.Range("Y901").Value = 3.75
$ws.Range("Z901").Value = -1
$ws.Range("AA901").Value = 0.95
$ws.Range("AB901").Value = 0.5125
$ws.Range("AC901").Value = -0.5

# Row 902
$ws.Range("B902").Value = 6838419
$ws.Range("F902").Value = 'Ajax'
$ws.Range("G902").Value = 'Heerenveen'
$ws.Range("H902").Value = 4
$ws.Range("I902").Value = 1
$ws.Range("J902").Value = 'H'
$ws.Range("K902").Value = 1.444
$ws.Range("L902").Value = 4.6
$ws.Range("M902").Value = 5.5
$ws.Range("N902").Value = 1.4
$ws.Range("O902").Value = 5.25
$ws.Range("P902").Value = 7
$ws.Range("Q902").Value = -1.5
$ws.Range("R902").Value = 2
$ws.Range("S902").Value = 1.85
$ws.Range("T902").Value = 3.5
$ws.Range("U902").Value = 1.975
$ws.Range("V902").Value = 1.875
$ws.Range("W902").Value = 0.3999999999999999
$ws.Range("Y902").Value = -1
$ws.Range("Z902").Value = 1
$ws.Range("AA902").Value = -1
$ws.Range("AB902").Value = 0.9750000000000001
$ws.Range("AC902").Value = -1

# Row 923
$ws.Range("B923").Value = 6838448
$ws.Range("F923").Value = 'RKC'
$ws.Range("G923").Value = 'Excelsior'
$ws.Range("H923").Value = 2
$ws.Range("I923").Value = 2
$ws.Range("J923").Value = 'D'
$ws.Range("K923").Value = 2.05
$ws.Range("L923").Value = 3.7
$ws.Range("M923").Value = 3.1
$ws.Range("N923").Value = 1.95
$ws.Range("O923").Value = 3.8
$ws.Range("P923").Value = 3.5
$ws.Range("R923").Value = 2
$ws.Range("S923").Value = 1.85
$ws.Range("T923").Value = 2.75
$ws.Range("U923").Value = 1.85
$ws.Range("V923").Value = 2
$ws.Range("W923").Value = -1
$ws.Range("X923").Value = 2.8
$ws.Range("Z923").Value = -1
$ws.Range("AA923").Value = 0.8500000000000001
$ws.Range("AB923").Value = 0.8500000000000001

# Row 924
$ws.Range("B924").Value = 6838447
$ws.Range("F924").Value = 'Fortuna Sittard'
$ws.Range("G924").Value = 'Vitesse'
$ws.Range("H924").Value = 3
$ws.Range("I924").Value = 1
$ws.Range("J924").Value = 'H'
$ws.Range("K924").Value = 2
$ws.Range("L924").Value = 3.5
$ws.Range("M924").Value = 3.4
$ws.Range("N924").Value = 2
$ws.Range("O924").Value = 3.6
$ws.Range("P924").Value = 3.6
$ws.Range("R924").Value = 2.025
$ws.Range("S924").Value = 1.825
$ws.Range("T924").Value = 2.5
$ws.Range("U924").Value = 1.875
$ws.Range("V924").Value = 1.975
$ws.Range("W924").Value = 1
$ws.Range("X924").Value = -1
$ws.Range("Z924").Value = 1.025
$ws.Range("AA924").Value = -1
$ws.Range("AB924").Value = 0.875

# Row 972
$ws.Range("B972").Value = 6838494
$ws.Range("F972").Value = 'Feyenoord'
$ws.Range("G972").Value = 'FC Twente'
$ws.Range("H972").Value = 0
$ws.Range("I972").Value = 0
$ws.Range("K972").Value = 1.666
$ws.Range("L972").Value = 3.9
$ws.Range("M972").Value = 4.75
$ws.Range("N972").Value = 1.45
$ws.Range("O972").Value = 4.5
$ws.Range("P972").Value = 6
$ws.Range("Q972").Value = -1.25
$ws.Range("R972").Value = 2.05
$ws.Range("S972").Value = 1.8
$ws.Range("U972").Value = 2.025
$ws.Range("V972").Value = 1.825
$ws.Range("X972").Value = 3.5
$ws.Range("Z972").Value = -1
$ws.Range("AA972").Value = 0.8
$ws.Range("AB972").Value = -1
$ws.Range("AC972").Value = 0.825

# Row 973
$ws.Range("B973").Value = 6838493
$ws.Range("F973").Value = 'Go Ahead Eagles'
$ws.Range("G973").Value = 'NEC'
$ws.Range("H973").Value = 2
$ws.Range("I973").Value = 2
$ws.Range("K973").Value = 2.4
$ws.Range("L973").Value = 3.5
$ws.Range("M973").Value = 2.75
$ws.Range("N973").Value = 2.2
$ws.Range("O973").Value = 3.5
$ws.Range("P973").Value = 3.2
$ws.Range("Q973").Value = -0.25
$ws.Range("R973").Value = 1.98
$ws.Range("S973").Value = 1.92
$ws.Range("U973").Value = 1.925
$ws.Range("V973").Value = 1.925
$ws.Range("X973").Value = 2.5
$ws.Range("Z973").Value = -0.5
$ws.Range("AA973").Value = 0.46
$ws.Range("AB973").Value = 0.925
$ws.Range("AC973").Value = -1

# Row 975
$ws.Range("R975").Value = 2.02
$ws.Range("S975").Value = 1.88
$ws.Range("T975").Value = 2.5
$ws.Range("U975").Value = 1.8
$ws.Range("V975").Value = 2.05

# Row 976
$ws.Range("R976").Value = 2.05
$ws.Range("S976").Value = 1.85

# Row 977
$ws.Range("N977").Value = 2.375
$ws.Range("P977").Value = 2.875
$ws.Range("R977").Value = 2.09
$ws.Range("S977").Value = 1.81
$ws.Range("U977").Value = 2
$ws.Range("V977").Value = 1.85

# Row 978
$ws.Range("R978").Value = 1.88
$ws.Range("S978").Value = 2.02

# Row 979
$ws.Range("R979").Value = 1.9
$ws.Range("S979").Value = 2

# Row 980
$ws.Range("R980").Value = 2.04
$ws.Range("S980").Value = 1.86

# Row 982
$ws.Range("N982").Value = 4.333
$ws.Range("P982").Value = 1.75
$ws.Range("R982").Value = 1.87
$ws.Range("S982").Value = 2.03

# Row 983
$ws.Range("N983").Value = 1.909
$ws.Range("P983").Value = 3.8
$ws.Range("R983").Value = 1.95
$ws.Range("S983").Value = 1.95
